$d = $word.ActiveDocument

# Locate the "-Shawn John" paragraph so we can insert the new
# "shawn2000jp@gmail.com" paragraph directly after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n", [char]7) -eq "-Shawn John") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '-Shawn John' paragraph"
}

# Build an insertion point just before the paragraph mark that ends the
# "-Shawn John" paragraph, so InsertXML() splits it into a new paragraph
# right after it (rather than merging with the following paragraph).
$insertPos = $target.Range.End - 1
$insertPoint = $d.Range($insertPos, $insertPos)

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:jc w:val="right"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:lang w:eastAsia="en-IN"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:lang w:eastAsia="en-IN"/>
    </w:rPr>
    <w:t>s</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:lang w:eastAsia="en-IN"/>
    </w:rPr>
    <w:t>hawn2000jp@gmail.com</w:t>
  </w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xmlFragment)
